$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (11) down into the
# three new rows so the new cells pick up the same date/time number
# formats (and no "All"-column style) as the rest of the table.
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E14").PasteSpecial(-4122)

# Row 12 - 28 Sep 2023, 10:00 - 12:00
$ws.Range("A12").Value = 45197
$ws.Range("B12").Value = 0.41666666666666669
$ws.Range("C12").Value = 0.5
$ws.Range("D12").Value = "All"
$ws.Range("E12").Value = "Discussion on the approach for weather station integration. Next meeting date scheduled."

# Row 13 - 2 Oct 2023, 13:00 - 14:00
$ws.Range("A13").Value = 45201
$ws.Range("B13").Value = 0.54166666666666663
$ws.Range("C13").Value = 0.58333333333333337
$ws.Range("D13").Value = "All"
$ws.Range("E13").Value = "summary and conclusion discussed and implemeneted"

# Row 14 - 5 Oct 2023, 21:00 - 23:00
$ws.Range("A14").Value = 45204
$ws.Range("B14").Value = 0.875
$ws.Range("C14").Value = 0.95833333333333337
$ws.Range("D14").Value = "All"
$ws.Range("E14").Value = "Report formatting done and finalised"

# View state: select E14 and scroll so row 3 is the top visible row
$ws.Range("E14").Select()
$excel.ActiveWindow.ScrollRow = 3
